# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - Row 3 ("R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 213
$wsOff.Range("C3").Value = 134
$wsOff.Range("D3").Value = 48
$wsOff.Range("E3").Value = 22
$wsOff.Range("F3").Value = 5

# DEF sheet - Row 3 ("R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 231
$wsDef.Range("C3").Value = 175
$wsDef.Range("D3").Value = 40
$wsDef.Range("E3").Value = 17
